$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Convert D2:D12 individual IF formulas into one shared formula (matches
#     what Excel does when the formula is re-entered/filled across the range) ---
$ws.Range("D2:D12").Formula = "=IF(C2>100,""Yes"",""No"")"

# --- Scientific-notation number format (same as the rest of column F / some D cells) ---
$sci = "0.00E+00"

# A couple of pre-existing Value cells picked up the scientific style too
$ws.Range("F6").NumberFormat = $sci
$ws.Range("F38").NumberFormat = $sci

# --- New data rows 39-43 (added in this commit) ---
# Row 39: Cu(H2O)(2,6-NDPA)0.5
$ws.Range("G39").Value = "10.1002/adfm.202007294"
$ws.Range("A39").Value = "Cu(H2O)(2,6-NDPA)0.5"
$ws.Range("B39").Value = "3D"
$ws.Range("D39").Value = "Yes"
$ws.Range("E39").Value = 2020
$ws.Range("F39").NumberFormat = $sci
$ws.Range("F39").Value = 2

# Row 40: Cu-BDC
$ws.Range("G40").Value = "10.1002/anie.201912642 "
$ws.Range("A40").Value = "Cu-BDC"
$ws.Range("B40").Value = "3D"
$ws.Range("C40").Value = 271
$ws.Range("D40").Value = "Yes"
$ws.Range("E40").Value = 2020
$ws.Range("F40").NumberFormat = $sci
$ws.Range("F40").Value = 0.01

# Row 41: HoHTTP
$ws.Range("G41").Value = "10.1038/s41557-019-0372-0 "
$ws.Range("A41").Value = "HoHTTP"
$ws.Range("B41").Value = "3D"
$ws.Range("C41").Value = 208
$ws.Range("D41").NumberFormat = $sci
$ws.Range("D41").Value = "Yes"
$ws.Range("E41").Value = 2020
$ws.Range("F41").NumberFormat = $sci
$ws.Range("F41").Value = 0.05

# Row 42: Cr(tri)2(CF3SO3)0.33
$ws.Range("G42").Value = "10.1038/s41557-021-00666-6 "
$ws.Range("A42").Value = "Cr(tri)2(CF3SO3)0.33"
$ws.Range("B42").Value = "3D"
$ws.Range("C42").Value = 80
$ws.Range("D42").NumberFormat = $sci
$ws.Range("D42").Value = "No"
$ws.Range("E42").Value = 2021
$ws.Range("F42").NumberFormat = $sci
$ws.Range("F42").Value = 0.014

# Row 43: Fe-HHTP
$ws.Range("A43").Value = "Fe-HHTP"
$ws.Range("G43").Value = "10.1002/anie.202102670 "
$ws.Range("B43").Value = "3D"
$ws.Range("C43").Value = 1400
$ws.Range("D43").NumberFormat = $sci
$ws.Range("D43").Value = "Yes"
$ws.Range("E43").Value = 2021
$ws.Range("F43").NumberFormat = $sci
$ws.Range("F43").Value = 0.001

# --- Final selection left on G47, matching the saved workbook state ---
$ws.Range("G47").Select() | Out-Null
